$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "57.604.34"
Set-TextValue "E2" "  -3.68%  "
Set-TextValue "D3" "2.275.59"
Set-TextValue "E3" "  -4.41%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "530.55"
Set-TextValue "E5" "  -4.45%  "
Set-TextValue "D6" "129.36"
Set-TextValue "E6" "  -2.81%  "
Set-TextValue "E8" "  -1.20%  "
Set-TextValue "D9" "2.269.10"
Set-TextValue "E9" "  -4.54%  "
Set-TextValue "D10" "0.0990"
Set-TextValue "E10" "  -5.39%  "
Set-TextValue "D11" "5.38"
Set-TextValue "E11" "  -4.42%  "
Set-TextValue "E12" "  -0.41%  "
Set-TextValue "E13" "  -4.48%  "
Set-TextValue "D14" "23.28"
Set-TextValue "E14" "  -4.47%  "
Set-TextValue "D15" "2.681.15"
Set-TextValue "E15" "  -4.58%  "
Set-TextValue "D16" "57.608.52"
Set-TextValue "E16" "  -3.61%  "
Set-TextValue "E17" "  -4.35%  "
Set-TextValue "D18" "2.274.75"
Set-TextValue "E18" "  -4.52%  "
Set-TextValue "D19" "10.43"
Set-TextValue "E19" "  -6.06%  "
Set-TextValue "D20" "4.19"
Set-TextValue "E20" "  -6.18%  "
Set-TextValue "D21" "312.13"
Set-TextValue "E21" "  -2.65%  "
Set-TextValue "D22" "6.31"
Set-TextValue "E22" "  -5.81%  "
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "62.31"
Set-TextValue "E24" "  -2.90%  "
Set-TextValue "E25" "  -4.60%  "
Set-TextValue "E26" "  -0.17%  "
Set-TextValue "D27" "7.96"
Set-TextValue "E27" "  -5.02%  "
Set-TextValue "D28" "1.28"
Set-TextValue "E28" "  -6.41%  "
Set-TextValue "D29" "170.29"
Set-TextValue "E29" "  +0.44%  "
Set-TextValue "E30" "  -5.13%  "
Set-TextValue "D31" "0.0₃0712"
Set-TextValue "E31" "  -5.94%  "
Set-TextValue "D32" "5.71"
Set-TextValue "E32" "  -5.33%  "
Set-TextValue "E33" "  -4.17%  "
Set-TextValue "B34" "USDe"
Set-TextValue "C34" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D34" "0.999"
Set-TextValue "E34" "  +0.01%  "
Set-TextValue "B35" "PolygonEcosystemToken"
Set-TextValue "C35" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D35" "0.374"
Set-TextValue "E35" "  -5.30%  "
Set-TextValue "D36" "17.68"
Set-TextValue "E36" "  -2.47%  "
Set-TextValue "E37" "  -0.16%  "
Set-TextValue "E38" "  -6.87%  "
Set-TextValue "D39" "3.86"
Set-TextValue "E39" "  -6.58%  "
Set-TextValue "D40" "37.92"
Set-TextValue "E41" "  -5.74%  "
Set-TextValue "D42" "139.94"
Set-TextValue "E42" "  -4.09%  "
Set-TextValue "D43" "285.09"
Set-TextValue "E43" "  -10.46%  "
Set-TextValue "E44" "  -3.76%  "
Set-TextValue "D45" "0.0943"
Set-TextValue "E45" "  -2.62%  "
Set-TextValue "E46" "  -3.05%  "
Set-TextValue "D47" "0.549"
Set-TextValue "E47" "  -3.83%  "
Set-TextValue "D48" "17.95"
Set-TextValue "E48" "  -8.83%  "
Set-TextValue "E49" "  -3.79%  "
Set-TextValue "D50" "10.94"
Set-TextValue "E50" "  -1.15%  "
Set-TextValue "E51" "  +84.43%  "

Write-Output "Applied all changes"
